$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 (header row for project 4): rename project and change version to numeric
$ws.Range("B20").Value = "Apache Maven Doxia"
$ws.Range("C20").Value = 1.6
$ws.Range("D20").Value = 13
$ws.Range("E20").Value = 51976
$ws.Range("F20").Value = 51.976

# Row 21
$ws.Range("C21").Value = 1.7
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = 52821
$ws.Range("F21").Value = 52.821

# Row 22
$ws.Range("C22").Value = 1.8
$ws.Range("D22").Value = 19
$ws.Range("E22").Value = 50236
$ws.Range("F22").Value = 50.236

# Row 23
$ws.Range("C23").Value = 1.9
$ws.Range("D23").Value = 12
$ws.Range("E23").Value = 58143
$ws.Range("F23").Value = 58.143

# Row 24: C24 keeps referencing the shared "version" label, which becomes "1.9.1"
$ws.Range("C24").Value = "1.9.1"
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = 58343
$ws.Range("F24").Value = 58.343
